# EIA Table 1.18.A: roll the report forward from "October 2016/2015" to
# "November 2016/2015" and update the revised data values (2017-01-31 update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_1_18_A")

# Helper: assign a literal text string to a cell/range without letting Excel's
# automatic data-type detection turn a "Month Year"-shaped string (e.g.
# "November 2016") into a date serial number. We stage the text in a distant
# scratch cell that has been explicitly formatted as Text ("@"), copy it, and
# paste-special "Values only" into the destination -- this preserves the
# destination's existing style/number format instead of creating a new one.
function Set-LiteralText {
    param(
        $TargetRange,
        [string]$Text
    )
    $helper = $ws.Range("Z100")
    $helper.NumberFormat = "@"
    $helper.Value = $Text
    $helper.Copy()
    $TargetRange.PasteSpecial(-4163)  # xlPasteValues
    $helper.Delete(-4159)             # xlShiftToLeft - remove helper cell entirely
}

# --- Report title / subtitle (row 2) ---
$ws.Range("A2").Value = "by State, by Sector, November 2016 and 2015 (Thousand Megawatthours)"

# --- Column period headers (row 6) ---
Set-LiteralText $ws.Range("B6") "November 2016"
Set-LiteralText $ws.Range("C6") "November 2015"
Set-LiteralText $ws.Range("E6") "November 2016"
Set-LiteralText $ws.Range("F6") "November 2015"
Set-LiteralText $ws.Range("G6") "November 2016"
Set-LiteralText $ws.Range("H6") "November 2015"
Set-LiteralText $ws.Range("I6") "November 2016"
Set-LiteralText $ws.Range("J6") "November 2015"
Set-LiteralText $ws.Range("K6") "November 2016"
Set-LiteralText $ws.Range("L6") "November 2015"

# --- Updated data values ---

# Row 32
$ws.Range("B32").Value = 4
$ws.Range("C32").Value = 4
$ws.Range("D32").Value = 0.14
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 4

# Row 35
$ws.Range("B35").Value = 4
$ws.Range("C35").Value = 4
$ws.Range("D35").Value = 0.14
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 4

# Row 52
$ws.Range("B52").Value = 40
$ws.Range("C52").Value = 54
$ws.Range("D52").Value = -0.255
$ws.Range("G52").Value = 40
$ws.Range("H52").Value = 54

# Row 53
$ws.Range("B53").Value = 35
$ws.Range("C53").Value = 46
$ws.Range("D53").Value = -0.242
$ws.Range("G53").Value = 35
$ws.Range("H53").Value = 46

# Row 57
$ws.Range("B57").Value = 5
$ws.Range("C57").Value = 7
$ws.Range("D57").Value = -0.34
$ws.Range("G57").Value = 5
$ws.Range("H57").Value = 7

# Row 61
$ws.Range("B61").Value = 140
$ws.Range("C61").Value = 147
$ws.Range("D61").Value = -0.048
$ws.Range("G61").Value = 140
$ws.Range("H61").Value = 147

# Row 62
$ws.Range("B62").Value = 140
$ws.Range("C62").Value = 147
$ws.Range("D62").Value = -0.048
$ws.Range("G62").Value = 140
$ws.Range("H62").Value = 147

# Row 68
$ws.Range("B68").Value = 184
$ws.Range("C68").Value = 204
$ws.Range("D68").Value = -0.099
$ws.Range("E68").Value = 4
$ws.Range("F68").Value = 4
$ws.Range("G68").Value = 180
$ws.Range("H68").Value = 201
